$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("D Green")
$ws.Range("B2").Value = 20
$ws.Range("B3").Value = 15
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = 0
$ws.Range("B7").ClearContents()
$ws.Range("B8").ClearContents()
$ws.Range("B9").Value = 0
$ws.Range("B10").Value = 13
$ws.Range("B11").Value = 0
$ws.Range("B12").ClearContents()
$ws.Range("B13").Value = 0
$ws.Range("B14").Value = 18
$ws.Range("B15").Value = 0

$ws = $wb.Worksheets.Item("Green")
$ws.Range("B2").Value = 32
$ws.Range("B3").Value = 70
$ws.Range("B4").Value = 53
$ws.Range("B5").Value = 29
$ws.Range("B6").Value = 37
$ws.Range("B7").Value = 96
$ws.Range("B8").Value = 32
$ws.Range("B9").Value = 152
$ws.Range("B10").Value = 60
$ws.Range("B11").Value = 44
$ws.Range("B12").Value = 32
$ws.Range("B13").Value = 60
$ws.Range("B14").Value = 51
$ws.Range("B15").Value = 100

$ws = $wb.Worksheets.Item("Yellow")
$ws.Range("B2").Value = 65.25
$ws.Range("B3").Value = 37.5
$ws.Range("B4").Value = 81.75
$ws.Range("B5").Value = 75.75
$ws.Range("B6").Value = 108.75
$ws.Range("B7").Value = 26.25
$ws.Range("B8").Value = 32.25
$ws.Range("B9").Value = 32.25
$ws.Range("B10").Value = 51
$ws.Range("B11").Value = 67.5
$ws.Range("B12").Value = 32.25
$ws.Range("B13").Value = 72
$ws.Range("B14").Value = 75
$ws.Range("B15").Value = 53.25

$ws = $wb.Worksheets.Item("Orange")
$ws.Range("B2").Value = 23
$ws.Range("B3").Value = 15.5
$ws.Range("B4").Value = 9
$ws.Range("B5").Value = 22
$ws.Range("B6").Value = 10
$ws.Range("B7").Value = 10
$ws.Range("B8").Value = 26
$ws.Range("B9").Value = 8
$ws.Range("B10").Value = 16.5
$ws.Range("B11").Value = 17.5
$ws.Range("B12").Value = 26
$ws.Range("B13").Value = 10
$ws.Range("B14").Value = 13.5
$ws.Range("B15").Value = 15.5

$ws = $wb.Worksheets.Item("Brown")
$ws.Range("B2").Value = 3.75
$ws.Range("B3").Value = 8.5
$ws.Range("B4").Value = 5.5
$ws.Range("B5").Value = 5
$ws.Range("B6").Value = 2.5
$ws.Range("B7").Value = 2
$ws.Range("B8").Value = 4.25
$ws.Range("B9").Value = 2.5
$ws.Range("B10").Value = 1
$ws.Range("B11").Value = 4
$ws.Range("B12").Value = 4.25
$ws.Range("B13").Value = 5.5
$ws.Range("B14").Value = 3.75
$ws.Range("B15").Value = 3.5

$ws = $wb.Worksheets.Item("Red")
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = 0
$ws.Range("B7").Value = 0
$ws.Range("B8").Value = 0
$ws.Range("B9").Value = 0
$ws.Range("B10").Value = 0
$ws.Range("B11").Value = 0
$ws.Range("B12").Value = 0
$ws.Range("B13").Value = 0
$ws.Range("B14").Value = 0
$ws.Range("B15").Value = 0

$ws = $wb.Worksheets.Item("Default Red")
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = 0
$ws.Range("B7").Value = 0
$ws.Range("B8").Value = 0
$ws.Range("B9").Value = 0
$ws.Range("B10").Value = 0
$ws.Range("B11").Value = 0
$ws.Range("B12").Value = 0
$ws.Range("B13").Value = 0
$ws.Range("B14").Value = 0
$ws.Range("B15").Value = 0

$ws = $wb.Worksheets.Item("Blue")
$ws.Range("B2").ClearContents()
$ws.Range("B3").ClearContents()
$ws.Range("B4").ClearContents()
$ws.Range("B5").ClearContents()
$ws.Range("B6").ClearContents()
$ws.Range("B7").ClearContents()
$ws.Range("B8").ClearContents()
$ws.Range("B9").ClearContents()
$ws.Range("B10").ClearContents()
$ws.Range("B11").ClearContents()
$ws.Range("B12").ClearContents()
$ws.Range("B13").ClearContents()
$ws.Range("B14").ClearContents()
$ws.Range("B15").ClearContents()

